# Weekly fruit/vegetable price update: insert 6 new daily records for
# 2022-01-18 (serial 44579) at row 109, shifting the existing rows 109:206
# down to 115:212 (dimension grows from A1:R206 to A1:R212).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows at row 109 (pushes old 109:206 -> 115:212, keeps
# column formatting such as the date style on column D).
$ws.Rows.Item(109).Resize(6).Insert()

# New row data: Mercado/Region/Categoria columns are constant for this
# sheet; only Variedad, Calidad, Volumen and the price columns vary.
$newRows = @(
    @{ Row=109; H="Calameño"; I="Extra";   J=5000; K=1500; L=1500; M=1500; P=1500 },
    @{ Row=110; H="Calameño"; I="Primera"; J=5000; K=1200; L=1200; M=1200; P=1200 },
    @{ Row=111; H="Calameño"; I="Segunda"; J=5000; K=1000; L=1000; M=1000; P=1000 },
    @{ Row=112; H="Tuna";     I="Extra";   J=5000; K=1500; L=1500; M=1500; P=1500 },
    @{ Row=113; H="Tuna";     I="Primera"; J=5000; K=1200; L=1200; M=1200; P=1200 },
    @{ Row=114; H="Tuna";     I="Segunda"; J=5000; K=1000; L=1000; M=1000; P=1000 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 4
    $ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
    $ws.Cells.Item($row, 3).Value = "Los Lagos"
    $ws.Cells.Item($row, 4).Value = (Get-Date -Year 2022 -Month 1 -Day 18 -Hour 0 -Minute 0 -Second 0)
    $ws.Cells.Item($row, 5).Value = 10
    $ws.Cells.Item($row, 6).Value = 100112027
    $ws.Cells.Item($row, 7).Value = "Melón"
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = "`$/unidad"
    $ws.Cells.Item($row, 15).Value = "Región de O'Higgins"
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = 1
    $ws.Cells.Item($row, 18).Value = "Hortaliza"
}
